{"js": "// Replace the arithmetic answers in the single 20x5 table with the new\n// values, cell by cell, preserving all existing paragraph/run formatting.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"values\");\nawait context.sync();\n\nconst newValues = [\n  [\"91-44=47\", \"52-28=24\", \"3+9=12\", \"5+89=94\", \"6+28=34\"],\n  [\"27+69=96\", \"45-8=37\", \"8+57=65\", \"56-39=17\", \"19+27=46\"],\n  [\"38+9=47\", \"54-6=48\", \"91-37=54\", \"29+16=45\", \"15+16=31\"],\n  [\"7+46=53\", \"30-9=21\", \"8+83=91\", \"43-24=19\", \"19+14=33\"],\n  [\"29+69=98\", \"46-17=29\", \"95-88=7\", \"23+29=52\", \"86-78=8\"],\n  [\"45+29=74\", \"26+38=64\", \"80-18=62\", \"63+28=91\", \"36+7=43\"],\n  [\"14-5=9\", \"94-65=29\", \"93-54=39\", \"80-35=45\", \"44+9=53\"],\n  [\"56+29=85\", \"7+69=76\", \"45+19=64\", \"44-38=6\", \"52-48=4\"],\n  [\"91-46=45\", \"70-23=47\", \"35+26=61\", \"74-35=39\", \"38+58=96\"],\n  [\"36+17=53\", \"23+29=52\", \"6+35=41\", \"6+89=95\", \"8+76=84\"],\n  [\"90-82=8\", \"87-38=49\", \"3+38=41\", \"39+27=66\", \"55-36=19\"],\n  [\"29+59=88\", \"9+44=53\", \"53+29=82\", \"36+45=81\", \"51-34=17\"],\n  [\"73-5=68\", \"96-77=19\", \"33-27=6\", \"61-48=13\", \"81-13=68\"],\n  [\"9+35=44\", \"48+48=96\", \"76-38=38\", \"33-25=8\", \"6+29=35\"],\n  [\"16+29=45\", \"90-1=89\", \"82-25=57\", \"3+29=32\", \"51-18=33\"],\n  [\"62-53=9\", \"20-16=4\", \"67+17=84\", \"37+6=43\", \"78+3=81\"],\n  [\"71-26=45\", \"9+57=66\", \"80-41=39\", \"17+66=83\", \"91-24=67\"],\n  [\"28+8=36\", \"64-38=26\", \"50-39=11\", \"61-29=32\", \"28+16=44\"],\n  [\"62-39=23\", \"41-9=32\", \"23-5=18\", \"6+57=63\", \"18+3=21\"],\n  [\"22-18=4\", \"14-8=6\", \"23+8=31\", \"70-69=1\", \"59+32=91\"],\n];\n\n// Sanity check the shape before writing so we fail loudly instead of\n// silently corrupting the table if the document doesn't match what we\n// expect.\nconst current = table.values;\nif (current.length !== newValues.length) {\n  throw new Error(\n    \"Unexpected row count: \" + current.length + \" vs \" + newValues.length\n  );\n}\nfor (let r = 0; r < current.length; r++) {\n  if (current[r].length !== newValues[r].length) {\n    throw new Error(\"Unexpected column count on row \" + r);\n  }\n}\n\ntable.values = newValues;\nawait context.sync();\n", "ps1": "# Replace the arithmetic answers in the single 20x5 table with the new\n# values, cell by cell, preserving all existing paragraph/run formatting.\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$newValues = @(\n    @(\"91-44=47\", \"52-28=24\", \"3+9=12\", \"5+89=94\", \"6+28=34\"),\n    @(\"27+69=96\", \"45-8=37\", \"8+57=65\", \"56-39=17\", \"19+27=46\"),\n    @(\"38+9=47\", \"54-6=48\", \"91-37=54\", \"29+16=45\", \"15+16=31\"),\n    @(\"7+46=53\", \"30-9=21\", \"8+83=91\", \"43-24=19\", \"19+14=33\"),\n    @(\"29+69=98\", \"46-17=29\", \"95-88=7\", \"23+29=52\", \"86-78=8\"),\n    @(\"45+29=74\", \"26+38=64\", \"80-18=62\", \"63+28=91\", \"36+7=43\"),\n    @(\"14-5=9\", \"94-65=29\", \"93-54=39\", \"80-35=45\", \"44+9=53\"),\n    @(\"56+29=85\", \"7+69=76\", \"45+19=64\", \"44-38=6\", \"52-48=4\"),\n    @(\"91-46=45\", \"70-23=47\", \"35+26=61\", \"74-35=39\", \"38+58=96\"),\n    @(\"36+17=53\", \"23+29=52\", \"6+35=41\", \"6+89=95\", \"8+76=84\"),\n    @(\"90-82=8\", \"87-38=49\", \"3+38=41\", \"39+27=66\", \"55-36=19\"),\n    @(\"29+59=88\", \"9+44=53\", \"53+29=82\", \"36+45=81\", \"51-34=17\"),\n    @(\"73-5=68\", \"96-77=19\", \"33-27=6\", \"61-48=13\", \"81-13=68\"),\n    @(\"9+35=44\", \"48+48=96\", \"76-38=38\", \"33-25=8\", \"6+29=35\"),\n    @(\"16+29=45\", \"90-1=89\", \"82-25=57\", \"3+29=32\", \"51-18=33\"),\n    @(\"62-53=9\", \"20-16=4\", \"67+17=84\", \"37+6=43\", \"78+3=81\"),\n    @(\"71-26=45\", \"9+57=66\", \"80-41=39\", \"17+66=83\", \"91-24=67\"),\n    @(\"28+8=36\", \"64-38=26\", \"50-39=11\", \"61-29=32\", \"28+16=44\"),\n    @(\"62-39=23\", \"41-9=32\", \"23-5=18\", \"6+57=63\", \"18+3=21\"),\n    @(\"22-18=4\", \"14-8=6\", \"23+8=31\", \"70-69=1\", \"59+32=91\")\n)\n\nfor ($r = 1; $r -le $t.Rows.Count; $r++) {\n    for ($c = 1; $c -le $t.Columns.Count; $c++) {\n        $cell = $t.Cell($r, $c)\n        $cell.Range.Text = $newValues[$r - 1][$c - 1]\n    }\n}\n"}
